$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point drift on the existing last row (A20)
$ws.Cells.Item(20, 1).Value = 45865.87529501157

# Append the new row of data (row 21)
$ws.Cells.Item(21, 1).Value = 45865.958641329
$ws.Cells.Item(21, 1).NumberFormat = $ws.Cells.Item(20, 1).NumberFormat

$ws.Cells.Item(21, 2).Value = 2025
$ws.Cells.Item(21, 3).Value = 30
$ws.Cells.Item(21, 4).Value = 13.89
$ws.Cells.Item(21, 5).Value = 91.16
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 3.09
$ws.Cells.Item(21, 8).Value = "N"
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = "23:00:26"
